# Automatische test-sync: 2025-06-29 15:13:50
# Append a new test-mail log entry (row 30) to the "Logs" sheet and
# bump the matching "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 30

$logs.Cells.Item($row, 1).Value = "Ik wil mijn gegevens aanpassen"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #15: Ik wil mijn gegevens aanpassen"
$logs.Cells.Item($row, 4).Value = "Overig"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nDank u voor uw bericht. Om uw gegevens aan te passen, kunt u inloggen op onze website met uw gebruikersnaam en wachtwoord. Eenmaal ingelogd, kunt u uw persoonlijke gegevens bijwerken onder uw accountinstellingen.`nMocht u nog verdere assistentie nodig hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-29 15:13:14"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"

# The multi-line answer in column E would otherwise leave the new row with
# an explicit "autofit-on-entry" height; re-autofit it so it collapses back
# to the sheet's standard (non-custom) row height, like the other rows.
$logs.Rows($row).AutoFit()

# Extend the existing conditional-formatting rules (D/G/H/I) so they keep
# covering the newly added row, without touching their rules/priorities.
$logs.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D30"))
$logs.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G30"))
$logs.Range("H2:H29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H30"))
$logs.Range("I2:I29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I30"))

# Update the Dashboard "Overig" tally (row 6, column B) to include the new entry.
$dashboard.Cells.Item(6, 2).Value = 3
